$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.99000000000062
$ws.Range("G2").Value = 0.000001973236989583071
$ws.Range("H2").Value = 0.0001242694811054434
$ws.Range("K2").Value = 5.499258905705647
$ws.Range("L2").Value = "[3.1697112300561683, 7.828806581355125]"
$ws.Range("M2").Value = 0.000004537565152951117
$ws.Range("N2").Value = 0.000009075130305902235
$ws.Range("O2").Value = -1.572368695490387
$ws.Range("P2").Value = "[-2.1006845771751568, -1.0440528138056173]"
$ws.Range("Q2").Value = 0.000000009239380149850263
$ws.Range("R2").Value = 0.00000001847876029970053
$ws.Range("S2").Value = 13.5146221760694
$ws.Range("T2").Value = "[12.047871575821958, 14.981372776316842]"
$ws.Range("W2").Value = 6.504004004004162
$ws.Range("X2").Value = 4.318658658658766
$ws.Range("Y2").Value = 8.689349349349559

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 23.07000000000017
$ws.Range("G3").Value = 0.0006577930309357161
$ws.Range("H3").Value = 0.003167988793610525
$ws.Range("K3").Value = 6.766257859292183
$ws.Range("L3").Value = "[2.3135395928933047, 11.21897612569106]"
$ws.Range("M3").Value = 0.003087627183227859
$ws.Range("N3").Value = 0.003087627183227859
$ws.Range("O3").Value = -1.962316131972003
$ws.Range("P3").Value = "[-2.6793162571156195, -1.2453160068283857]"
$ws.Range("Q3").Value = 0.0000001998096494926216
$ws.Range("R3").Value = 0.0000001998096494926216
$ws.Range("S3").Value = 15.04226245106362
$ws.Range("T3").Value = "[12.641670714862414, 17.44285418726482]"
$ws.Range("W3").Value = 7.205045045045097
$ws.Range("X3").Value = 4.572432432432464
$ws.Range("Y3").Value = 9.837657657657729
